# Updates the crypto price list (coin/link/price/volume columns) in place,
# matching the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.860.93"
$ws.Range("E2").Value = "  -5.16%  "
$ws.Range("D3").Value = "3.321.65"
$ws.Range("E3").Value = "  -6.20%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.02"
$ws.Range("E5").Value = "  -4.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.62"
$ws.Range("E6").Value = "  -7.69%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("E8").Value = "  -3.98%  "
$ws.Range("D9").Value = "3.311.96"
$ws.Range("E9").Value = "  -6.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("E10").Value = "  -10.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.581"
$ws.Range("E11").Value = "  -7.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.97"
$ws.Range("E12").Value = "  -9.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -8.20%  "
$ws.Range("D14").Value = "3.860.76"
$ws.Range("E14").Value = "  -5.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.51"
$ws.Range("E15").Value = "  -7.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "597.86"
$ws.Range("E16").Value = "  -10.29%  "
$ws.Range("D17").Value = "66.006.81"
$ws.Range("E17").Value = "  -5.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.85"
$ws.Range("E18").Value = "  -3.27%  "
$ws.Range("D19").Value = "3.334.21"
$ws.Range("E19").Value = "  -6.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.116"
$ws.Range("E20").Value = "  -4.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.31"
$ws.Range("E21").Value = "  -9.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.896"
$ws.Range("E22").Value = "  -7.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.73"
$ws.Range("E23").Value = "  -7.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.05"
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.52"
$ws.Range("E25").Value = "  -4.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.99"
$ws.Range("E26").Value = "  -8.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.99"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.63"
$ws.Range("E28").Value = "  -9.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -9.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.62"
$ws.Range("E30").Value = "  -10.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.45"
$ws.Range("E31").Value = "  -8.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.19"
$ws.Range("E32").Value = "  -8.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.71"
$ws.Range("E33").Value = "  -14.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.94"
$ws.Range("E34").Value = "  -6.88%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.103"
$ws.Range("E35").Value = "  -6.65%  "
$ws.Range("D36").Value = "3.771.28"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "57.71"
$ws.Range("E37").Value = "  -6.62%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "529.58"
$ws.Range("E39").Value = "  +5.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.42"
$ws.Range("E40").Value = "  -7.29%  "
$ws.Range("D41").Value = "0.0₃0706"
$ws.Range("E41").Value = "  -13.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.64"
$ws.Range("E42").Value = "  -9.02%  "
$ws.Range("E43").Value = "  -7.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.337"
$ws.Range("E44").Value = "  -8.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.53"
$ws.Range("E45").Value = "  -9.03%  "
$ws.Range("E46").Value = "  -3.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0410"
$ws.Range("E47").Value = "  -8.40%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.128"
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("E49").Value = "  -9.46%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.89"
$ws.Range("E51").Value = "  +7.12%  "
